$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "conversion" column M: (E<row> * 15625) / 60
# Entered as two fill operations (M13:M15 and M17:M20) with M16/M21 typed
# individually, mirroring how the values were produced, then the whole new
# range is reset back to the workbook's default (unformatted) style.
$ws.Range("M13:M15").Formula = "=(E13*15625)/60"
$ws.Range("M16").Formula = "=(E16*15625)/60"
$ws.Range("M17:M20").Formula = "=(E17*15625)/60"
$ws.Range("M21").Formula = "=(E21*15625)/60"

$ws.Range("M13:M21").Style = "Normal"

# Move the active selection as it ended up after the edit.
$ws.Range("L9").Select() | Out-Null
